$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first fixture row (Real Madrid CF - Athletic Club de Bilbao / 17/10/2021)
# This shifts all subsequent rows up by one.
$ws.Rows.Item(1).Delete()

# Append the new fixture at the end (now row 10)
$ws.Range("A10").Value = "Real Madrid CF - Granada CF"

# "06/02/2022" looks like a date, so force text entry to avoid Excel's
# automatic date conversion, then restore the default cell style so the
# cell formatting matches the rest of the sheet.
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "06/02/2022"
$ws.Range("B10").Style = "Normal"
